# Tag info now pulled from a database file instead of being manually entered
# by the user: remove the manually-entered "Name", "Description" and
# "Normal Price" columns (B:D) and add a sample data row sourced from the
# database.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Name", "Description" and "Normal Price" columns entirely.
# This shifts "Sale Price" (and everything after it) three columns to the
# left, carrying its number-format style along with it.
$ws.Range("B1:D1").EntireColumn.Delete()

# Populate the new data row pulled from the database.
$ws.Range("A2").Value = 1001000
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 7

# Reset the view/selection back to the top-left of the sheet.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
